$wb = $excel.ActiveWorkbook

# Sheet "展览" - update "想去人数" (interested count) for three events
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 3149
$ws1.Range("F4").Value = 148
$ws1.Range("F5").Value = 121

# Sheet "全部类型" - same events duplicated, keep them in sync
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F7").Value = 3149
$ws4.Range("F8").Value = 148
$ws4.Range("F10").Value = 121
